$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee")

# Fill in July / September monthly scorecard figures (column C = July,
# column D = August already had data, column E = September) for each
# tracked metric row. August (D) values are unchanged but re-written so the
# alignment normalises below.
$ws.Cells.Item(28, 3).Value = 21
$ws.Cells.Item(28, 4).Value = 22
$ws.Cells.Item(28, 5).Value = 19

$ws.Cells.Item(29, 3).Value = 11
$ws.Cells.Item(29, 4).Value = 22
$ws.Cells.Item(29, 5).Value = 19

$ws.Cells.Item(30, 3).Value = 10
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0

$ws.Cells.Item(31, 3).Value = 6
$ws.Cells.Item(31, 4).Value = 16
$ws.Cells.Item(31, 5).Value = 3

$ws.Cells.Item(32, 3).Value = 5
$ws.Cells.Item(32, 4).Value = 6
$ws.Cells.Item(32, 5).Value = 9

$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(33, 5).Value = 0

$ws.Cells.Item(34, 3).Value = 3
$ws.Cells.Item(34, 4).Value = 2
$ws.Cells.Item(34, 5).Value = 9

$ws.Cells.Item(35, 3).Value = 2
$ws.Cells.Item(35, 4).Value = 4
$ws.Cells.Item(35, 5).Value = 0

# Centre-align the whole data block (C:E now match each other visually).
$ws.Range("C28:E35").HorizontalAlignment = -4108

# Leave the cursor where the author last left it.
$ws.Range("E32").Select()
